$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.462.55'
$ws.Range('E2').Value = '  -1.68%  '
$ws.Range('D3').Value = '2.420.90'
$ws.Range('E3').Value = '  -2.22%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.30'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '143.09'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '0.529'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').Value = '2.420.36'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  -5.03%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.20'
$ws.Range('E12').Value = '  -2.61%  '
$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  -3.32%  '
$ws.Range('D14').Value = '26.56'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').Value = '0.0000173'
$ws.Range('E15').Value = '  -6.12%  '
$ws.Range('D16').Value = '2.868.59'
$ws.Range('E16').Value = '  -3.46%  '
$ws.Range('D17').Value = '62.369.30'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').Value = '2.427.87'
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').Value = '11.12'
$ws.Range('E19').Value = '  -4.19%  '
$ws.Range('D20').Value = '7.19'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('D21').Value = '324.49'
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').Value = '2.06'
$ws.Range('E23').Value = '  +8.14%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '65.14'
$ws.Range('E25').Value = '  -3.73%  '
$ws.Range('D26').Value = '611.72'
$ws.Range('E26').Value = '  -3.84%  '
$ws.Range('D27').Value = '8.88'
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').Value = '0.0₃0976'
$ws.Range('E28').Value = '  -7.32%  '
$ws.Range('D29').Value = '2.553.02'
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').Value = '1.46'
$ws.Range('E31').Value = '  -3.87%  '
$ws.Range('D32').Value = '8.05'
$ws.Range('E32').Value = '  -5.05%  '
$ws.Range('D33').Value = '1.87'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('E34').Value = '  -5.30%  '
$ws.Range('D35').Value = '5.04'
$ws.Range('E35').Value = '  -3.50%  '
$ws.Range('D36').Value = '1.47'
$ws.Range('E36').Value = '  -5.37%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '0.374'
$ws.Range('E38').Value = '  -3.28%  '
$ws.Range('D39').Value = '18.69'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '147.19'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').Value = '5.24'
$ws.Range('E41').Value = '  -5.12%  '
$ws.Range('D42').Value = '1.74'
$ws.Range('E42').Value = '  -6.17%  '
$ws.Range('D43').Value = '2.54'
$ws.Range('E43').Value = '  -3.75%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '42.11'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').Value = '144.57'
$ws.Range('E46').Value = '  -4.22%  '
$ws.Range('D47').Value = '3.71'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.30'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('D49').Value = '0.0525'
$ws.Range('E49').Value = '  -4.83%  '
$ws.Range('D50').Value = '0.594'
$ws.Range('E50').Value = '  -2.64%  '
$ws.Range('D51').Value = '0.0228'
$ws.Range('E51').Value = '  -5.10%  '
